$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-07-23 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-24 Thursday", 2) | Out-Null

# Update the multiplication problems in the table, cell by cell, since some
# source values are duplicated within the same row (e.g. row 10) and a global
# Find/Replace would not be able to target them independently.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "76×65="
$t.Cell(1, 2).Range.Text = "57×65="
$t.Cell(1, 3).Range.Text = "40×65="
$t.Cell(1, 4).Range.Text = "51×53="
$t.Cell(1, 5).Range.Text = "47×77="

$t.Cell(5, 1).Range.Text = "88×99="
$t.Cell(5, 2).Range.Text = "71×76="
$t.Cell(5, 3).Range.Text = "20×84="
$t.Cell(5, 4).Range.Text = "81×43="
$t.Cell(5, 5).Range.Text = "97×97="

$t.Cell(10, 1).Range.Text = "91×33="
$t.Cell(10, 2).Range.Text = "48×29="
$t.Cell(10, 3).Range.Text = "22×34="
$t.Cell(10, 4).Range.Text = "62×50="
$t.Cell(10, 5).Range.Text = "13×83="

$t.Cell(15, 1).Range.Text = "97×44="
$t.Cell(15, 2).Range.Text = "87×63="
$t.Cell(15, 3).Range.Text = "58×12="
$t.Cell(15, 4).Range.Text = "46×35="
$t.Cell(15, 5).Range.Text = "35×14="

$t.Cell(20, 1).Range.Text = "85×60="
$t.Cell(20, 2).Range.Text = "94×85="
$t.Cell(20, 3).Range.Text = "76×30="
$t.Cell(20, 4).Range.Text = "22×26="
$t.Cell(20, 5).Range.Text = "23×20="

Write-Output "edits applied"
